$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1814345991561181
$ws.Range("C2").Value = 0.5780590717299579
$ws.Range("J2").Value = 0.01687763713080169
$ws.Range("P2").Value = 0.1350210970464135
$ws.Range("S2").Value = 0.08860759493670886
$ws.Range("B3").Value = 0.01438848920863309
$ws.Range("C3").Value = 0.01438848920863309
$ws.Range("J3").Value = 0.007194244604316547
$ws.Range("P3").Value = 0.7266187050359713
$ws.Range("S3").Value = 0.237410071942446
$ws.Range("J4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.7045454545454546
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("B6").Value = 0.06578947368421052
$ws.Range("D6").Value = 0.01754385964912281
$ws.Range("F6").Value = 0.07017543859649122
$ws.Range("J6").Value = 0.2456140350877193
$ws.Range("O6").Value = 0.01754385964912281
$ws.Range("Q6").Value = 0.1535087719298246
$ws.Range("R6").Value = 0.1008771929824561
$ws.Range("S6").Value = 0.3289473684210527
$ws.Range("B7").Value = 0.1036269430051813
$ws.Range("D7").Value = 0.0155440414507772
$ws.Range("E7").Value = 0.005181347150259068
$ws.Range("F7").Value = 0.04663212435233161
$ws.Range("J7").Value = 0.150259067357513
$ws.Range("O7").Value = 0.005181347150259068
$ws.Range("Q7").Value = 0.1917098445595855
$ws.Range("R7").Value = 0.04145077720207254
$ws.Range("S7").Value = 0.4404145077720207
$ws.Range("B8").Value = 0.06986899563318777
$ws.Range("D8").Value = 0.0240174672489083
$ws.Range("F8").Value = 0.06768558951965066
$ws.Range("J8").Value = 0.1026200873362445
$ws.Range("O8").Value = 0.01965065502183406
$ws.Range("Q8").Value = 0.1790393013100437
$ws.Range("R8").Value = 0.09606986899563319
$ws.Range("S8").Value = 0.4410480349344978
$ws.Range("B9").Value = 0.09142857142857143
$ws.Range("D9").Value = 0.01714285714285714
$ws.Range("F9").Value = 0.04571428571428571
$ws.Range("J9").Value = 0.09142857142857143
$ws.Range("O9").Value = 0.02857142857142857
$ws.Range("Q9").Value = 0.1657142857142857
$ws.Range("R9").Value = 0.08571428571428572
$ws.Range("S9").Value = 0.4742857142857143
$ws.Range("B10").Value = 0.08527755430410297
$ws.Range("D10").Value = 0.02172164119066774
$ws.Range("E10").Value = 0.0008045052292839903
$ws.Range("F10").Value = 0.08125502815768303
$ws.Range("J10").Value = 0.1053901850362027
$ws.Range("O10").Value = 0.01448109412711183
$ws.Range("Q10").Value = 0.2244569589702333
$ws.Range("R10").Value = 0.08608205953338696
$ws.Range("S10").Value = 0.3805309734513274
$ws.Range("G11").Value = 0.1234939759036145
$ws.Range("J11").Value = 0.105421686746988
$ws.Range("K11").Value = 0.2018072289156627
$ws.Range("L11").Value = 0.5632530120481928
$ws.Range("S11").Value = 0.006024096385542169
$ws.Range("G12").Value = 0.6989795918367347
$ws.Range("J12").Value = 0.1836734693877551
$ws.Range("K12").Value = 0.00510204081632653
$ws.Range("L12").Value = 0.03571428571428571
$ws.Range("S12").Value = 0.07653061224489796
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.375
$ws.Range("F15").Value = 0.01945525291828794
$ws.Range("H15").Value = 0.178988326848249
$ws.Range("I15").Value = 0.09727626459143969
$ws.Range("J15").Value = 0.3346303501945525
$ws.Range("K15").Value = 0.08171206225680934
$ws.Range("M15").Value = 0.01167315175097276
$ws.Range("O15").Value = 0.07392996108949416
$ws.Range("S15").Value = 0.2023346303501946
$ws.Range("F16").Value = 0.0189873417721519
$ws.Range("H16").Value = 0.1708860759493671
$ws.Range("I16").Value = 0.06329113924050633
$ws.Range("J16").Value = 0.379746835443038
$ws.Range("K16").Value = 0.1645569620253164
$ws.Range("M16").Value = 0.006329113924050633
$ws.Range("N16").Value = 0.006329113924050633
$ws.Range("O16").Value = 0.08227848101265822
$ws.Range("S16").Value = 0.1075949367088608
$ws.Range("F17").Value = 0.02603036876355748
$ws.Range("H17").Value = 0.2082429501084599
$ws.Range("I17").Value = 0.08459869848156182
$ws.Range("J17").Value = 0.3796095444685466
$ws.Range("K17").Value = 0.08676789587852494
$ws.Range("M17").Value = 0.01301518438177874
$ws.Range("O17").Value = 0.07592190889370933
$ws.Range("S17").Value = 0.1258134490238612
$ws.Range("F18").Value = 0.0101010101010101
$ws.Range("H18").Value = 0.1919191919191919
$ws.Range("I18").Value = 0.08585858585858586
$ws.Range("J18").Value = 0.4343434343434344
$ws.Range("K18").Value = 0.09090909090909091
$ws.Range("M18").Value = 0.0101010101010101
$ws.Range("O18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.1212121212121212
$ws.Range("F19").Value = 0.01778808971384378
$ws.Range("H19").Value = 0.1956689868522815
$ws.Range("I19").Value = 0.06496519721577726
$ws.Range("J19").Value = 0.3735498839907193
$ws.Range("K19").Value = 0.1206496519721578
$ws.Range("M19").Value = 0.0185614849187935
$ws.Range("N19").Value = 0.002320185614849188
$ws.Range("O19").Value = 0.08275328692962104
$ws.Range("S19").Value = 0.1237432327919567
